$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.375.79'
$ws.Range('E2').Value = '  +1.42%  '

$ws.Range('D3').Value = '1.828.52'
$ws.Range('E3').Value = '  +0.39%  '

$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').Value = '314.07'
$ws.Range('E5').Value = '  +1.02%  '

$ws.Range('E6').Value = '  +0.05%  '

$ws.Range('D7').Value = '0.4472'
$ws.Range('E7').Value = '  +5.19%  '

$ws.Range('D8').Value = '0.3767'
$ws.Range('E8').Value = '  +3.08%  '

$ws.Range('D9').Value = '0.07538'
$ws.Range('E9').Value = '  +4.41%  '

$ws.Range('D10').Value = '0.8952'
$ws.Range('E10').Value = '  +6.49%  '

$ws.Range('D11').Value = '21.02'
$ws.Range('E11').Value = '  +2.34%  '

$ws.Range('D12').Value = '1.812.64'
$ws.Range('E12').Value = '  -0.41%  '

$ws.Range('D13').Value = '6.752'
$ws.Range('E13').Value = '  +1.72%  '

$ws.Range('D14').Value = '94.55'
$ws.Range('E14').Value = '  +5.81%  '

$ws.Range('D15').Value = '5.401'
$ws.Range('E15').Value = '  +2.46%  '

$ws.Range('D16').Value = '0.07118'
$ws.Range('E16').Value = '  +0.95%  '

$ws.Range('E17').Value = '  +0.01%  '

$ws.Range('D18').Value = '0.000008828'
$ws.Range('E18').Value = '  +1.05%  '

$ws.Range('D20').Value = '15.22'
$ws.Range('E20').Value = '  +2.68%  '

$ws.Range('D21').Value = '27.387.95'
$ws.Range('E21').Value = '  +0.96%  '

$ws.Range('D22').Value = '5.286'
$ws.Range('E22').Value = '  +3.36%  '

$ws.Range('D23').Value = '10.96'
$ws.Range('E23').Value = '  +1.58%  '

$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').Value = '2.004'
$ws.Range('E24').Value = '  +1.43%  '

$ws.Range('B25').Value = 'LidoDAOToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D25').Value = '2.502'
$ws.Range('E25').Value = '  +12.68%  '

$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '151.56'
$ws.Range('E26').Value = '  +0.43%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '18.64'
$ws.Range('E27').Value = '  +2.69%  '

$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = '5.365'
$ws.Range('E28').Value = '  +3.02%  '

$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').Value = '118.04'
$ws.Range('E29').Value = '  +1.35%  '

$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').Value = '0.08852'
$ws.Range('E30').Value = '  +1.58%  '

$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = '0.7812'
$ws.Range('E31').Value = '  +6.61%  '

$ws.Range('B32').Value = 'ARBITRUM'
$ws.Range('C32').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D32').Value = '1.203'
$ws.Range('E32').Value = '  +2.55%  '

$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '4.477'
$ws.Range('E33').Value = '  +1.55%  '

$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = '2.893'
$ws.Range('E34').Value = '  -0.30%  '

$ws.Range('B35').Value = 'Frax'
$ws.Range('C35').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D35').Value = '1.001'
$ws.Range('E35').Value = '  +0.04%  '

$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').Value = '1.108'
$ws.Range('E36').Value = '  +1.86%  '

$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '0.01989'
$ws.Range('E37').Value = '  +2.76%  '

$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = '0.05332'
$ws.Range('E38').Value = '  +2.39%  '

$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').Value = '7.391'
$ws.Range('E39').Value = '  +2.50%  '

$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '0.5331'
$ws.Range('E40').Value = '  +4.61%  '

$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').Value = '0.1729'
$ws.Range('E41').Value = '  +2.82%  '

$ws.Range('D42').Value = '2.873'
$ws.Range('E42').Value = '  +0.37%  '

$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = '2.273'
$ws.Range('E43').Value = '  +16.47%  '

$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').Value = '8.826'
$ws.Range('E44').Value = '  +3.62%  '

$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = '0.5163'
$ws.Range('E45').Value = '  +9.42%  '

$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '10.78'
$ws.Range('E46').Value = '  +2.29%  '

$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').Value = '106.33'
$ws.Range('E47').Value = '  +0.58%  '

$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = '1.707'
$ws.Range('E48').Value = '  +3.79%  '

$ws.Range('B49').Value = 'PaxDollar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D49').Value = '1.001'
$ws.Range('E49').Value = '  +0.09%  '

$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.06379'
$ws.Range('E50').Value = '  +1.08%  '

$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '64.52'
$ws.Range('E51').Value = '  +3.49%  '
